$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1) Title paragraph: merge "OPBM " + "Pre-requisite and Installer Warnings"
#    into a single run with the same combined text (find/replace auto-merges
#    same-format runs into one run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("OPBM Pre-requisite and Installer Warnings", $true, $false, $false, $false, $false,
                         $true, 1, $false, "OPBM Pre-requisite and Installer Warnings", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Date paragraph: "August 16, 2011" -> three runs: "August ", "18", ", 2011"
#    Use InsertXML on the paragraph's text range (excluding the paragraph
#    mark) so the exact run boundaries from the target are produced.
# ---------------------------------------------------------------------------
$dateParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "August 16, 2011") {
        $dateParaIndex = $i
        break
    }
}
$pDate = $d.Paragraphs.Item($dateParaIndex)
$startDate = $pDate.Range.Start
$endDate = $pDate.Range.End - 1
$rDate = $d.Range($startDate, $endDate)
$rDate.InsertXML("<w:p $wns><w:r><w:t xml:space='preserve'>August </w:t></w:r><w:r><w:t>18</w:t></w:r><w:r><w:t>, 2011</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 3) Remove the _GoBack bookmark from the "Rick C. Hodgin" paragraph (it no
#    longer starts that paragraph).
# ---------------------------------------------------------------------------
$rickIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Rick C. Hodgin") {
        $rickIndex = $i
        break
    }
}
$pRick = $d.Paragraphs.Item($rickIndex)
$pRick.Range.InsertXML("<w:p $wns><w:r><w:t>Rick C. Hodgin</w:t></w:r></w:p>")

# ---------------------------------------------------------------------------
# 4) Java Virtual Machine paragraph: merge the 5 runs into a single run.
# ---------------------------------------------------------------------------
$javaOld = "Java Virtual Machine installed (JDK or JRE version 1.7.0 32-bit or 64-bit, preferred 64-bit)"
$d.Content.Find.Execute($javaOld, $true, $false, $false, $false, $false,
                         $true, 1, $false, $javaOld, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Internet Explorer paragraph: drop the leading space.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" Internet Explorer 9 (will not work with IE8 or older)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Internet Explorer 9 (will not work with IE8 or older)", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Insert a new "24-bit video resolution of at least 800 x 600" list item
#    right after "Microsoft Office 2010" (inherits its numPr automatically).
# ---------------------------------------------------------------------------
$msOfficeIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Microsoft Office 2010") {
        $msOfficeIndex = $i
        break
    }
}
$pMsOffice = $d.Paragraphs.Item($msOfficeIndex)
$pMsOffice.Range.InsertParagraphAfter()
$pVideo = $d.Paragraphs.Item($msOfficeIndex + 1)
$pVideo.Range.Text = "24-bit video resolution of at least 800 x 600"

# ---------------------------------------------------------------------------
# 7) Move the _GoBack bookmark into the second blank "ListParagraph ind=0"
#    paragraph that follows (the one right before "OPBM will install:").
# ---------------------------------------------------------------------------
$installIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "OPBM will install:") {
        $installIndex = $i
        break
    }
}
$pBlankBookmark = $d.Paragraphs.Item($installIndex - 1)
$pBlankBookmark.Range.InsertXML("<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='0'/></w:pPr><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>")

# ---------------------------------------------------------------------------
# 8) Insert a new "7-zip" list item right before "Adobe Acrobat Reader X
#    (version 10)" (inherits its numPr automatically).
# ---------------------------------------------------------------------------
$adobeIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq "Adobe Acrobat Reader X (version 10)") {
        $adobeIndex = $i
        break
    }
}
$pAdobe = $d.Paragraphs.Item($adobeIndex)
$pAdobe.Range.InsertParagraphBefore()
$pZip = $d.Paragraphs.Item($adobeIndex)
$pZip.Range.Text = "7-zip"

Write-Host "Edit complete."
